$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 27; $row++) {
    # Column H: PERIOD TO EXPIRE -- decreases by 1 (one more day has elapsed)
    $hCell = $ws.Cells.Item($row, 8)
    $hCell.Value2 = $hCell.Value2 - 1

    # Column I: LAST UPDATE -- bump the progress date to 04-Nov-2025.
    # A direct string assignment ("04-Nov-2025") gets auto-parsed by Excel
    # as a real date (changing the cell's type/format), so instead write it
    # as a text formula and then paste-special as a value to freeze it back
    # down to a plain text literal without disturbing the cell's style.
    $iCell = $ws.Cells.Item($row, 9)
    $iCell.Formula = '="04-Nov-2025"'
    $iCell.Copy() | Out-Null
    $iCell.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

$excel.CutCopyMode = $false
